$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet was renamed from "DisabilityTypes Data" to "Data"
$ws.Name = "Data"

# The schema changed: the table now only holds the three "Инвалид ... группы"
# rows (Id 3/4/5) instead of the original ten disability-type rows, and the
# EndDate sentinel moved from 9999-12-31 to 9998-12-31 (serial 2958465 ->
# 2958100). Drop the now-unused rows 5-11 first.
$ws.Range("A5:D11").Clear()

# Row 2 -> Id 3 / "Инвалид I группы"
$ws.Range("A2").NumberFormat = "@"
$ws.Cells.Item(2, 1).Value = "3"
$ws.Cells.Item(2, 2).Value = "Инвалид I группы"
$ws.Cells.Item(2, 3).Value = 2
$ws.Cells.Item(2, 4).Value = 2958100

# Row 3 -> Id 4 / "Инвалид II группы"
$ws.Range("A3").NumberFormat = "@"
$ws.Cells.Item(3, 1).Value = "4"
$ws.Cells.Item(3, 2).Value = "Инвалид II группы"
$ws.Cells.Item(3, 3).Value = 2
$ws.Cells.Item(3, 4).Value = 2958100

# Row 4 -> Id 5 / "Инвалид III группы"
$ws.Range("A4").NumberFormat = "@"
$ws.Cells.Item(4, 1).Value = "5"
$ws.Cells.Item(4, 2).Value = "Инвалид III группы"
$ws.Cells.Item(4, 3).Value = 2
$ws.Cells.Item(4, 4).Value = 2958100

# Column A held text-typed ids ("3","4","5") in the original file; drop the
# temporary text number-format now that the values are entered so the cells
# go back to their normal (unstyled) look, same as before.
$ws.Range("A2:A4").ClearFormats()
